$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy date-column style (s="2": bold, thin border, centered, custom date numfmt) from an
# existing, unaffected cell (A92) onto the new cells that need it: the inserted row 93
# and the newly appended rows 113-115.
$ws.Range("A92").Copy()
$ws.Range("A93").PasteSpecial(-4122)
$ws.Range("A113:A115").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A90").Value = 44232
$ws.Range("B90").Value = 1
$ws.Range("C90").Value = 11
$ws.Range("D90").Value = 101.3357899585444

$ws.Range("A91").Value = 44233
$ws.Range("B91").Value = 1
$ws.Range("C91").Value = 13
$ws.Range("D91").Value = 119.7604790419162

$ws.Range("A92").Value = 44234
$ws.Range("B92").Value = 5
$ws.Range("C92").Value = 14
$ws.Range("D92").Value = 128.972823583602

$ws.Range("A93").Value = 44235
$ws.Range("B93").Value = 1
$ws.Range("C93").Value = 14
$ws.Range("D93").Value = 128.972823583602

$ws.Range("A94").Value = 44236
$ws.Range("B94").Value = 3
$ws.Range("C94").Value = 19
$ws.Range("D94").Value = 175.0345462920313

$ws.Range("A95").Value = 44237
$ws.Range("B95").Value = 2
$ws.Range("C95").Value = 20
$ws.Range("D95").Value = 184.2468908337172

$ws.Range("A96").Value = 44238
$ws.Range("B96").Value = 1
$ws.Range("C96").Value = 15
$ws.Range("D96").Value = 138.1851681252879

$ws.Range("A97").Value = 44239
$ws.Range("B97").Value = 6
$ws.Range("C97").Value = 17
$ws.Range("D97").Value = 156.6098572086596

$ws.Range("A98").Value = 44240
$ws.Range("B98").Value = 2
$ws.Range("C98").Value = 16
$ws.Range("D98").Value = 147.3975126669737

$ws.Range("A99").Value = 44241
$ws.Range("B99").Value = 0
$ws.Range("C99").Value = 14
$ws.Range("D99").Value = 128.972823583602

$ws.Range("A100").Value = 44242
$ws.Range("B100").Value = 3
$ws.Range("C100").Value = 14
$ws.Range("D100").Value = 128.972823583602

$ws.Range("A101").Value = 44243
$ws.Range("B101").Value = 2
$ws.Range("C101").Value = 11
$ws.Range("D101").Value = 101.3357899585444

$ws.Range("A102").Value = 44244
$ws.Range("B102").Value = 0
$ws.Range("C102").Value = 10
$ws.Range("D102").Value = 92.12344541685859

$ws.Range("A103").Value = 44245
$ws.Range("B103").Value = 1
$ws.Range("C103").Value = 14
$ws.Range("D103").Value = 128.972823583602

$ws.Range("A104").Value = 44246
$ws.Range("B104").Value = 3
$ws.Range("C104").Value = 14
$ws.Range("D104").Value = 128.972823583602

$ws.Range("A105").Value = 44247
$ws.Range("B105").Value = 1
$ws.Range("C105").Value = 14
$ws.Range("D105").Value = 128.972823583602

$ws.Range("A106").Value = 44248
$ws.Range("B106").Value = 4
$ws.Range("C106").Value = 14
$ws.Range("D106").Value = 128.972823583602

$ws.Range("A107").Value = 44249
$ws.Range("B107").Value = 3
$ws.Range("C107").Value = 18
$ws.Range("D107").Value = 165.8222017503455

$ws.Range("A108").Value = 44250
$ws.Range("B108").Value = 2
$ws.Range("C108").Value = 20
$ws.Range("D108").Value = 184.2468908337172

$ws.Range("A109").Value = 44251
$ws.Range("B109").Value = 0
$ws.Range("C109").Value = 22
$ws.Range("D109").Value = 202.6715799170889

$ws.Range("A110").Value = 44252
$ws.Range("B110").Value = 5
$ws.Range("C110").Value = 20
$ws.Range("D110").Value = 184.2468908337172

$ws.Range("A111").Value = 44253
$ws.Range("B111").Value = 5
$ws.Range("C111").Value = 24
$ws.Range("D111").Value = 221.0962690004606

$ws.Range("A112").Value = 44254
$ws.Range("B112").Value = 3
$ws.Range("C112").Value = 31
$ws.Range("D112").Value = 285.5826807922617

$ws.Range("A113").Value = 44255
$ws.Range("B113").Value = 2

$ws.Range("A114").Value = 44256
$ws.Range("B114").Value = 7

$ws.Range("A115").Value = 44257
$ws.Range("B115").Value = 9
